# Applies the edits described by the commit diff:
#  - rank abbreviation change
#  - surname run split update (Byś -> Bys)
#  - date shifts (2021-04-18 -> 2021-05-11, 2021-04-19 -> 2021-05-12)
#  - city change (Dębica -> Warszawa)
#  - fill-in of the (previously empty) transport-means run

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

$d.Content.Find.Execute("szer. pchor.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "kpr. pchor.", $wdReplaceAll)

$d.Content.Find.Execute("Paweł Byś", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Paweł", $wdReplaceAll)

$d.Content.Find.Execute(" Byś", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, " Bys", $wdReplaceAll)

$d.Content.Find.Execute("2021-04-18", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "2021-05-11", $wdReplaceAll)

$d.Content.Find.Execute("2021-04-19", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "2021-05-12", $wdReplaceAll)

$d.Content.Find.Execute("Dębica", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Warszawa", $wdReplaceAll)

$d.Content.Find.Execute("zbiorowym ()", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "zbiorowym (kolejowym w klasie 2, w pociągu ekspresowym)", $wdReplaceAll)
